# Update cryptocurrency price/volume data cells to reflect refreshed market figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.665.85"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "2.239.28"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "2.581.11"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "2.330.86"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.824"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "44.433.05"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "0.0₃0926"
$ws.Range("E19").Value = "  -4.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0773"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0295"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "1.807.14"
$ws.Range("E44").Value = "  +3.53%  "
$ws.Range("E45").Value = "  +8.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "81.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.61%  "
